$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the word count (B2); dependent formulas (B4, D9, E9, D11, E11) recalc automatically
$ws.Range("B2").Value = 19358

# Update the saved selection/active cell shown in the sheet view
$ws.Range("G21").Select()
